$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.330.37"
$ws.Range("E2").Value = "  -0.71%  "

# Row 3
$ws.Range("D3").Value = "3.319.70"
$ws.Range("E3").Value = "  -1.50%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "187.80"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.06%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "559.56"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.34%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.12%  "

# Row 8
$ws.Range("E8").Value = "  -1.59%  "

# Row 9
$ws.Range("D9").Value = "3.312.19"
$ws.Range("E9").Value = "  -1.50%  "

# Row 10
$ws.Range("E10").Value = "  -1.43%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.588"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.40%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.74"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.54%  "

# Row 13
$ws.Range("E13").Value = "  +1.31%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.67"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.05%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "633.86"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.64%  "

# Row 16
$ws.Range("D16").Value = "3.850.27"
$ws.Range("E16").Value = "  -1.39%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.15"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.50%  "

# Row 18
$ws.Range("D18").Value = "66.267.44"
$ws.Range("E18").Value = "  -0.62%  "

# Row 19
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.117"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.90%  "

# Row 20
$ws.Range("B20").Value = "WrappedEther"
$ws.Range("C20").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D20").Value = "3.276.81"
$ws.Range("E20").Value = "  -2.90%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.02"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -6.00%  "

# Row 22
$ws.Range("E22").Value = "  -0.45%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.33"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +7.57%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.14"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +7.21%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.97"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.61%  "

# Row 26
$ws.Range("E26").Value = "  -3.44%  "

# Row 27
$ws.Range("E27").Value = "  +0.63%  "

# Row 28
$ws.Range("E28").Value = "  -0.84%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.63"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.31%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.71"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.42%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.42"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.57%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.98"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.58%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.42"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.75%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.11"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.42%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "555.12"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.14%  "

# Row 36
$ws.Range("E36").Value = "  -0.45%  "

# Row 37
$ws.Range("D37").Value = "3.839.48"
$ws.Range("E37").Value = "  +1.19%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "57.76"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.69%  "

# Row 39
$ws.Range("E39").Value = "  +0.04%  "

# Row 40
$ws.Range("E40").Value = "  +0.89%  "

# Row 41
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "34.01"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.54%  "

# Row 42
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.31"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.74%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.72"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.21%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.129"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.16%  "

# Row 45
$ws.Range("B45").Value = "CoreDAO"
$ws.Range("C45").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.28"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -13.43%  "

# Row 46
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.337"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.64%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0422"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.82%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.24"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.57%  "

# Row 49
$ws.Range("E49").Value = "  -0.93%  "

# Row 50
$ws.Range("E50").Value = "  -2.92%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.03%  "
